# Add a new row to the "Change History" sheet of the PPC0003 template
# describing CRE21-011: "To revise the sampling criteria for post
# payment check report PPC0003", dated 2021/08/03.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change History")

$newRow = 8

# Seed the new row's cells by copying the formatting (and, for column A,
# the numeric typing) from the most similar existing rows, mirroring how
# this template's previous Change History entries were authored.
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("B7").Copy($ws.Range("B8"))
$ws.Range("C5").Copy($ws.Range("C8"))
$ws.Range("D7").Copy($ws.Range("D8"))

# Fill in the new Change History entry (Item 5 / CRE21-011).
$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "CRE21-011"
$ws.Cells.Item($newRow, 3).Value = "To revise the sampling criteria for post payment check report PPC0003"
$ws.Cells.Item($newRow, 4).Value = "2021/08/03"
